# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$updates = @{
    7  = -3
    9  = 4
    12 = 1
    15 = 4
    25 = -2
    32 = 3
    34 = 1
    39 = 2
    53 = 2
    54 = 1
    66 = -1
    72 = -2
    74 = 3
    76 = -1
    77 = -2
    80 = -2
    86 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
